# Test Data Updated to match passwords correctly
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (yuvi.subu1@yopmail.com) had a mismatched password "Welcome@12";
# correct it to the real password "Welcome@123" used by the other rows.
$ws.Range("B3").Value = "Welcome@123"
